# LoginCredentialDetails.xlsx -- loginPage test changes / irrigation module data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Update TestCaseId / credential values for rows 3-7 (new "SU-T7xx" cases)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "SU-T721"
$ws.Range("B3").Value = "superadminD"
$ws.Range("C3").Value = "Rihand@123"

$ws.Range("A4").Value = "SU-T722"
$ws.Range("C4").Value = "Rihand@123"

$ws.Range("A5").Value = "SU-T723"
$ws.Range("B5").Value = "superadminD"

$ws.Range("A6").Value = "SU-T730"
$ws.Range("B6").ClearContents()

$ws.Range("A7").Value = "SU-T731"
$ws.Range("C7").ClearContents()

# New Email column entries for rows 6 & 7
$ws.Range("D6").Value = "awtindia.sc@gmail.com"
$ws.Range("D7").Value = "niclosTesla@gmail.com"

# ---------------------------------------------------------------------------
# 2. Clear out the old trailing rows (8-12) content but keep their formatting,
#    then add a fresh blank formatted row 13
# ---------------------------------------------------------------------------
$ws.Range("A8").ClearContents()
$ws.Range("D8").ClearContents()

$ws.Range("A9").ClearContents()
$ws.Range("D9").ClearContents()

$ws.Range("A10").ClearContents()
$ws.Range("C10").ClearContents()

$ws.Range("A11").ClearContents()
$ws.Range("C11").ClearContents()

$ws.Range("A12").ClearContents()
$ws.Range("C12").ClearContents()

$ws.Range("A13").Value = ""
$ws.Range("C13").Value = ""

# ---------------------------------------------------------------------------
# 3. Fix up cell formatting that changed alongside the content
# ---------------------------------------------------------------------------
# B3 / C3 lose their special (Admin / Hyperlink) styling -> plain Normal
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Style = "Normal"

# D6/D7/D8 (and the new A13/C13 cells) pick up the small style (xf index 4)
# already used by C6/C7/C10-C12 -- copy formats across.
$ws.Range("C6").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D8").PasteSpecial(-4122)

$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("C13").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Rebuild the hyperlinks: drop the old ones, add the new set that matches
#    the refreshed credential / email cells
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:Rihand@123")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:Testing@123")
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:niclosTesla@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:awtindia.sc@gmail.com")

# Adding a hyperlink re-applies the Hyperlink cell style (and mutates the
# font), so restore the original (non-hyperlink-flavoured) styles that these
# cells carried before -- copy the format back in from sibling cells that
# already have the right xf.
$ws.Range("C5").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C5").PasteSpecial(-4122)

$ws.Range("C6").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D7").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5. Sheet view / selection: scroll back to top-left and select D10
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D10").Select()
